$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New body text for the "e099a" test case (retreating into enemy territory / a battle)
$bodyText = @'
<Bold>e099a Retreat into a Battle</Bold> 
<InlineUIContainer><Button Content='r11.33' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
You retreated into a battle.
<LineBreak/><LineBreak/>
                        <InlineUIContainer><Image Name='Sherman1' Height='200' Width='325'></Image></InlineUIContainer>
<LineBreak/><LineBreak/>
Click image to continue.  
'@

# Insert a new row right above the "e100" entry (row 109) to hold the new
# "e099a" test case, shifting the existing rows (e100 onward) down by one.
$ws.Rows(109).Insert()

$ws.Range("A109").Value = "e099a"
$ws.Range("B109").Value = $bodyText

# Match the row height used for similarly-sized entries (120pt).
$ws.Rows(109).RowHeight = 120

# Reflect the new selection position used when the author saved the file.
[void]$ws.Range("B110").Select()
